$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 276
$ws.Range("F7").Value = 13180
$ws.Range("F8").Value = 68
$ws.Range("F10").Value = 291
$ws.Range("F11").Value = 4395
$ws.Range("F12").Value = 6787
$ws.Range("F13").Value = 63
$ws.Range("F15").Value = 3538
$ws.Range("F16").Value = 43
$ws.Range("F17").Value = 11
$ws.Range("F18").Value = 16
$ws.Range("F19").Value = 177
$ws.Range("F21").Value = 43
$ws.Range("F23").Value = 130
$ws.Range("F24").Value = 3684
$ws.Range("F26").Value = 3902
$ws.Range("F27").Value = 3902
$ws.Range("F28").Value = 425
$ws.Range("F29").Value = 1929
$ws.Range("F31").Value = 253
$ws.Range("F32").Value = 6929
$ws.Range("F34").Value = 167
$ws.Range("F35").Value = 1961
$ws.Range("F36").Value = 2044
$ws.Range("F38").Value = 115
$ws.Range("F39").Value = 1090
$ws.Range("F41").Value = 230
$ws.Range("F43").Value = 226
$ws.Range("F46").Value = 6
$ws.Range("F47").Value = 148
$ws.Range("F48").Value = 1845
$ws.Range("F49").Value = 75
$ws.Range("F50").Value = 169

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 129

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 476
$ws.Range("F3").Value = 647
$ws.Range("F4").Value = 34

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 476
$ws.Range("F6").Value = 647
$ws.Range("F7").Value = 34
$ws.Range("F8").Value = 276
$ws.Range("F9").Value = 13180
$ws.Range("F10").Value = 68
$ws.Range("F13").Value = 291
$ws.Range("F14").Value = 4395
$ws.Range("F15").Value = 6787
$ws.Range("F16").Value = 63
$ws.Range("F17").Value = 3538
$ws.Range("F18").Value = 43
$ws.Range("F19").Value = 11
$ws.Range("F20").Value = 16
$ws.Range("F22").Value = 43
$ws.Range("F25").Value = 129
$ws.Range("F26").Value = 130
$ws.Range("F28").Value = 3902
$ws.Range("F29").Value = 425
$ws.Range("F31").Value = 253
$ws.Range("F32").Value = 6929
$ws.Range("F35").Value = 167
$ws.Range("F36").Value = 1962
$ws.Range("F37").Value = 2044
$ws.Range("F39").Value = 115
$ws.Range("F40").Value = 1090
$ws.Range("F41").Value = 230
$ws.Range("F42").Value = 226
$ws.Range("F45").Value = 148
$ws.Range("F47").Value = 1845
$ws.Range("F48").Value = 75
$ws.Range("F50").Value = 169
